$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster table (player, position, team) for rows 2-19
$data = @(
    @("Trae Young",        "PG",       "Atlanta Hawks"),
    @("Norman Powell",     "SG,SF",    "LA Clippers"),
    @("Jared McCain",      "PG,SG",    "Philadelphia 76ers"),
    @("LeBron James",      "SF,PF",    "Los Angeles Lakers"),
    @("Nicolas Claxton",   "C",        "Brooklyn Nets"),
    @("Jeremy Sochan",     "SF,PF",    "San Antonio Spurs"),
    @("Walker Kessler",    "C",        "Utah Jazz"),
    @("Dereck Lively II",  "C",        "Dallas Mavericks"),
    @("Alperen Sengün",    "C",        "Houston Rockets"),
    @("Shaedon Sharpe",    "SG,SF",    "Portland Trail Blazers"),
    @("Desmond Bane",      "SG,SF",    "Memphis Grizzlies"),
    @("Devin Booker",      "PG,SG",    "Phoenix Suns"),
    @("Jalen Brunson",     "PG",       "New York Knicks"),
    @("P.J. Washington",   "PF",       "Dallas Mavericks"),
    @("Devin Vassell",     "SG,SF",    "San Antonio Spurs"),
    @("Immanuel Quickley", "PG,SG",    "Toronto Raptors"),
    @("Kawhi Leonard",     "SG,SF,PF", "LA Clippers"),
    @("Coby White",        "PG,SG",    "Chicago Bulls")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
